# Auto-generated edit script applying cryptos.xlsx diff updates (Fri Oct 18 11:58:19 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.721.82'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '2.624.03'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''595.40'
$ws.Range("E5").Value = '  +0.61%  '
$ws.Range("D6").Value = '''153.09'
$ws.Range("E6").Value = '  +1.14%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.98%  '
$ws.Range("D9").Value = '2.623.47'
$ws.Range("E9").Value = '  +1.01%  '
$ws.Range("E10").Value = '  +9.77%  '
$ws.Range("E11").Value = '  -0.52%  '
$ws.Range("E12").Value = '  +1.93%  '
$ws.Range("E13").Value = '  +0.89%  '
$ws.Range("D14").Value = '''27.49'
$ws.Range("E14").Value = '  +0.82%  '
$ws.Range("E15").Value = '  +5.23%  '
$ws.Range("D16").Value = '3.100.74'
$ws.Range("E16").Value = '  +0.86%  '
$ws.Range("D17").Value = '67.691.29'
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("D18").Value = '2.616.71'
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("D19").Value = '''11.37'
$ws.Range("E19").Value = '  +3.82%  '
$ws.Range("D20").Value = '''366.18'
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("E21").Value = '  +1.21%  '
$ws.Range("E22").Value = '  -1.49%  '
$ws.Range("E23").Value = '  -0.72%  '
$ws.Range("E24").Value = '  +2.36%  '
$ws.Range("D25").Value = '''71.92'
$ws.Range("E25").Value = '  +8.67%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").Value = '''9.83'
$ws.Range("E27").Value = '  -0.51%  '
$ws.Range("E29").Value = '  +3.32%  '
$ws.Range("D30").Value = '''1.01'
$ws.Range("E30").Value = '  +0.77%  '
$ws.Range("D31").Value = '''573.34'
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("D32").Value = '''7.90'
$ws.Range("E32").Value = '  +3.18%  '
$ws.Range("D33").Value = '''1.39'
$ws.Range("E33").Value = '  +2.10%  '
$ws.Range("E34").Value = '  +1.64%  '
$ws.Range("E35").Value = '  +5.38%  '
$ws.Range("D36").Value = '''1.00'
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").Value = '''1.55'
$ws.Range("E37").Value = '  +4.43%  '
$ws.Range("D38").Value = '''160.15'
$ws.Range("E38").Value = '  +2.79%  '
$ws.Range("D39").Value = '''19.09'
$ws.Range("E39").Value = '  +1.10%  '
$ws.Range("E40").Value = '  +5.39%  '
$ws.Range("D41").Value = '''0.366'
$ws.Range("E41").Value = '  +0.61%  '
$ws.Range("D42").Value = '''5.32'
$ws.Range("E42").Value = '  +2.59%  '
$ws.Range("B43").Value = 'BabyDogeCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D43").Value = '0.0₆0333'
$ws.Range("E43").Value = '  +17.32%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = '''2.64'
$ws.Range("E44").Value = '  +4.98%  '
$ws.Range("E45").Value = '  +5.41%  '
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").Value = '''40.09'
$ws.Range("E47").Value = '  -1.69%  '
$ws.Range("E48").Value = '  +0.84%  '
$ws.Range("E49").Value = '  -0.80%  '
$ws.Range("D50").Value = '''21.79'
$ws.Range("E50").Value = '  +2.66%  '
$ws.Range("E51").Value = '  +0.12%  '
